$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.665.65"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.881.47"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4820"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2837"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "1.994.50"
$ws.Range("E10").Value = "  +6.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07500"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6662"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "30.630.47"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.247.60"
$ws.Range("E18").Value = "  +5.45%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007616"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +16.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.307"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.198"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.334"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.947"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.432"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09573"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.343"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.049"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05037"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.218"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7490"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.708"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01859"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.094"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9160"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4288"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.822"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.447"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1293"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.479"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.967"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05651"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.66%  "
